$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# passportNumber (C2): plain text, no special style
$ws.Range("C2").Value = "211NWCIND"

# visaNumber (D2): text value that looks numeric, but must NOT carry the
# quote-prefix ("stored as text") style - enter with a leading apostrophe
# (forces text) then clear the resulting quote-prefix formatting.
$ws.Range("D2").Value = "'12929"
$ws.Range("D2").ClearFormats()

# date (J2): text value "7 " (trailing space) keeping quote-prefix style
$ws.Range("J2").Value = "'7 "

# timeSlot (K2): text value " 10:00" keeping quote-prefix style
$ws.Range("K2").Value = "' 10:00"

[void]$ws.Range("J14").Select()
